$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.896.59'

$ws.Range('E2').Value = '  +2.44%  '

$ws.Range('D3').Value = '3.411.09'

$ws.Range('E3').Value = '  +3.14%  '

$ws.Range('E4').Value = '  -0.01%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '577.29'
$cell.Style = "Normal"

$ws.Range('E5').Value = '  +2.80%  '

$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '138.34'
$cell.Style = "Normal"

$ws.Range('E6').Value = '  +7.29%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '3.411.74'

$ws.Range('E8').Value = '  +3.13%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.478'
$cell.Style = "Normal"

$ws.Range('E9').Value = '  +0.95%  '

$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '7.50'
$cell.Style = "Normal"

$ws.Range('E10').Value = '  +2.07%  '

$ws.Range('E11').Value = '  +9.28%  '

$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '0.396'
$cell.Style = "Normal"

$ws.Range('E12').Value = '  +6.86%  '

$ws.Range('D13').Value = '3.990.53'

$ws.Range('E13').Value = '  +3.61%  '

$ws.Range('E14').Value = '  +1.81%  '

$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.0000181'
$cell.Style = "Normal"

$ws.Range('E15').Value = '  +7.95%  '

$ws.Range('D16').Value = '3.409.60'

$ws.Range('E16').Value = '  +3.64%  '

$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '25.45'
$cell.Style = "Normal"

$ws.Range('E17').Value = '  +4.86%  '

$ws.Range('D18').Value = '61.904.21'

$ws.Range('E18').Value = '  +2.29%  '

$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '14.14'
$cell.Style = "Normal"

$ws.Range('E19').Value = '  +6.04%  '

$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '5.92'
$cell.Style = "Normal"

$ws.Range('E20').Value = '  +4.44%  '

$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '9.45'
$cell.Style = "Normal"

$ws.Range('E21').Value = '  +5.15%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '391.80'
$cell.Style = "Normal"

$ws.Range('E22').Value = '  +11.61%  '

$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '0.573'
$cell.Style = "Normal"

$ws.Range('E23').Value = '  +3.38%  '

$ws.Range('D24').Value = '3.545.67'

$ws.Range('E24').Value = '  +3.45%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '0.0000129'
$cell.Style = "Normal"

$ws.Range('E25').Value = '  +18.89%  '

$ws.Range('E26').Value = '  +0.24%  '

$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '71.47'
$cell.Style = "Normal"

$ws.Range('E27').Value = '  +3.02%  '

$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '1.60'
$cell.Style = "Normal"

$ws.Range('E28').Value = '  +11.03%  '

$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '7.70'
$cell.Style = "Normal"

$ws.Range('E29').Value = '  +5.91%  '

$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"

$ws.Range('E30').Value = '  +0.25%  '

$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '8.30'
$cell.Style = "Normal"

$ws.Range('E31').Value = '  +5.71%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '0.159'
$cell.Style = "Normal"

$ws.Range('E32').Value = '  +5.38%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '2.17'
$cell.Style = "Normal"

$ws.Range('E33').Value = '  +3.22%  '

$ws.Range('B34').Value = 'RenzoRestakedETH'

$ws.Range('C34').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'

$ws.Range('D34').Value = '3.440.90'

$ws.Range('E34').Value = '  +3.29%  '

$ws.Range('B35').Value = 'USDe'

$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'

$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"

$ws.Range('E35').Value = '  -0.01%  '

$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '23.58'
$cell.Style = "Normal"

$ws.Range('E36').Value = '  +3.83%  '

$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '5.53'
$cell.Style = "Normal"

$ws.Range('E37').Value = '  +4.93%  '

$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '6.98'
$cell.Style = "Normal"

$ws.Range('E38').Value = '  +2.88%  '

$ws.Range('E39').Value = '  +4.69%  '

$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '161.45'
$cell.Style = "Normal"

$ws.Range('E40').Value = '  +2.13%  '

$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.0795'
$cell.Style = "Normal"

$ws.Range('E41').Value = '  +5.98%  '

$ws.Range('E42').Value = '  +12.80%  '

$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"

$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('B44').Value = 'ONDO'

$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'

$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '1.23'
$cell.Style = "Normal"

$ws.Range('E44').Value = '  +6.84%  '

$ws.Range('B45').Value = 'Mantle'

$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'

$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '0.775'
$cell.Style = "Normal"

$ws.Range('E45').Value = '  +4.60%  '

$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '4.48'
$cell.Style = "Normal"

$ws.Range('E46').Value = '  +2.69%  '

$ws.Range('B47').Value = 'EnergySwap'

$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '25.20'
$cell.Style = "Normal"

$ws.Range('E47').Value = '  +10.62%  '

$ws.Range('B48').Value = 'OKB'

$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'

$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '41.57'
$cell.Style = "Normal"

$ws.Range('E48').Value = '  +1.07%  '

$ws.Range('E49').Value = '  +4.33%  '

$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '22.91'
$cell.Style = "Normal"

$ws.Range('E50').Value = '  +5.72%  '

$ws.Range('D51').Value = '2.375.33'

$ws.Range('E51').Value = '  +9.48%  '
